$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'30.471.00"
$ws.Range("E2").Value = "  -0.20%  "
$ws.Range("D3").Value = "'1.900.54"
$ws.Range("E3").Value = "  +1.39%  "
$ws.Range("E4").Value = "  +0.15%  "
$ws.Range("D5").Value = "'238.93"
$ws.Range("E5").Value = "  +0.99%  "
$ws.Range("D6").Value = "'1.000"
$ws.Range("E6").Value = "  +0.16%  "
$ws.Range("D7").Value = "'0.4905"
$ws.Range("E7").Value = "  +0.72%  "
$ws.Range("E8").Value = "  +0.88%  "
$ws.Range("D9").Value = "'0.06678"
$ws.Range("E9").Value = "  +0.16%  "
$ws.Range("D10").Value = "'1.873.18"
$ws.Range("E10").Value = "  -0.07%  "
$ws.Range("D11").Value = "'16.95"
$ws.Range("E11").Value = "  +1.99%  "
$ws.Range("D12").Value = "'0.07331"
$ws.Range("E12").Value = "  +1.50%  "
$ws.Range("D13").Value = "'5.178"
$ws.Range("E13").Value = "  +3.43%  "
$ws.Range("D14").Value = "'88.06"
$ws.Range("E14").Value = "  -1.42%  "
$ws.Range("D15").Value = "'0.6678"
$ws.Range("E15").Value = "  +2.19%  "
$ws.Range("D16").Value = "'30.445.70"
$ws.Range("E16").Value = "  -0.11%  "
$ws.Range("D17").Value = "'0.000007873"
$ws.Range("D18").Value = "'13.42"
$ws.Range("E18").Value = "  +3.14%  "
$ws.Range("D19").Value = "'1.000"
$ws.Range("E19").Value = "  +0.15%  "
$ws.Range("D20").Value = "'5.429"
$ws.Range("E20").Value = "  +14.91%  "
$ws.Range("D21").Value = "'2.143.00"
$ws.Range("E21").Value = "  +1.41%  "
$ws.Range("E22").Value = "  +0.17%  "
$ws.Range("D23").Value = "'194.79"
$ws.Range("E23").Value = "  -8.86%  "
$ws.Range("D24").Value = "'6.133"
$ws.Range("E24").Value = "  +0.14%  "
$ws.Range("D25").Value = "'9.496"
$ws.Range("D26").Value = "'162.59"
$ws.Range("E26").Value = "  +4.21%  "
$ws.Range("D27").Value = "'18.37"
$ws.Range("E27").Value = "  -3.43%  "
$ws.Range("D28").Value = "'1.938"
$ws.Range("E28").Value = "  +6.06%  "
$ws.Range("D29").Value = "'1.487"
$ws.Range("E29").Value = "  +5.78%  "
$ws.Range("D30").Value = "'4.330"
$ws.Range("E30").Value = "  +1.55%  "
$ws.Range("D31").Value = "'0.09147"
$ws.Range("E31").Value = "  +1.27%  "
$ws.Range("D32").Value = "'4.119"
$ws.Range("E32").Value = "  +4.96%  "
$ws.Range("D33").Value = "'0.05159"
$ws.Range("E33").Value = "  +0.94%  "
$ws.Range("D34").Value = "'0.7383"
$ws.Range("E34").Value = "  +2.04%  "
$ws.Range("D35").Value = "'1.106"
$ws.Range("E35").Value = "  +2.62%  "
$ws.Range("E36").Value = "  +1.65%  "
$ws.Range("D37").Value = "'0.01850"
$ws.Range("E37").Value = "  +2.14%  "
$ws.Range("D38").Value = "'2.674"
$ws.Range("E38").Value = "  +0.57%  "
$ws.Range("D39").Value = "'0.9249"
$ws.Range("E39").Value = "  +0.63%  "
$ws.Range("D40").Value = "'2.070"
$ws.Range("E40").Value = "  +1.19%  "
$ws.Range("D41").Value = "'0.4406"
$ws.Range("E41").Value = "  +0.20%  "
$ws.Range("D42").Value = "'106.87"
$ws.Range("D43").Value = "'5.899"
$ws.Range("E43").Value = "  +2.84%  "
$ws.Range("D44").Value = "'0.9954"
$ws.Range("E44").Value = "  +0.15%  "
$ws.Range("D45").Value = "'69.03"
$ws.Range("E45").Value = "  +21.09%  "
$ws.Range("D46").Value = "'0.1371"
$ws.Range("E46").Value = "  +3.37%  "
$ws.Range("D47").Value = "'7.565"
$ws.Range("E47").Value = "  +3.12%  "
$ws.Range("D48").Value = "'9.011"
$ws.Range("E48").Value = "  +4.10%  "
$ws.Range("D49").Value = "'34.97"
$ws.Range("E49").Value = "  +5.58%  "
$ws.Range("D50").Value = "'0.05835"
$ws.Range("E50").Value = "  +0.19%  "
$ws.Range("D51").Value = "'0.3930"
$ws.Range("E51").Value = "  -2.15%  "
